$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The values for columns D, K, L, M, N, O, P, Q, R, S, T got cyclically
# rotated across rows 2, 3, 4, 5, 7:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 7
#   new row 5 <- old row 3
#   new row 7 <- old row 5

# Capture the "before" values first so we don't clobber data we still need to read.
$rows = @(2, 3, 4, 5, 7)
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# new_row -> source_row (old values to copy in)
$mapping = @{
    2 = 4
    3 = 2
    4 = 7
    5 = 3
    7 = 5
}

foreach ($newRow in $rows) {
    $srcRow = $mapping[$newRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $srcData[$c]
    }
}
